# Apply the TextBlob bias-dataset edit: fix the two typo'd header strings,
# add TP/FP/TN/FN classification columns (F:I) with per-row formulas,
# add COUNTIF summary formulas (K:N) in row 2, fix header labels, and
# update the sheet selection to match the authored state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Fix the two mangled shared strings used by D1/E1 (typos from an
# earlier edit: "1itive_result" / "-1ative_result" -> the intended
# "Positive_result" / "Negative_result"). ---
$ws.Range("D1").Value = "Positive_result"
$ws.Range("E1").Value = "Negative_result"

# --- New header row for the classification + summary columns. ---
$ws.Range("F1").Value = "TP"
$ws.Range("G1").Value = "FP"
$ws.Range("H1").Value = "TN"
$ws.Range("I1").Value = "FN"

$ws.Range("K1").Value = "TP"
$ws.Range("L1").Value = "FP"
$ws.Range("M1").Value = "TN"
$ws.Range("N1").Value = "FN"

# --- Per-row TP/FP/TN/FN classification formulas, rows 2..68. ---
for ($row = 2; $row -le 68; $row++) {
    $ws.Range("F$row").Formula = "=IF(AND(A$row=-1,C$row=-1),`"TP`")"
    $ws.Range("G$row").Formula = "=IF(AND(A$row=1,C$row=-1),`"FP`")"
    $ws.Range("H$row").Formula = "=IF(AND(A$row=1,C$row=1),`"TN`")"
    $ws.Range("I$row").Formula = "=IF(AND(A$row=-1,C$row=1),`"FN`")"
}

# --- Summary COUNTIF formulas, row 2 only. ---
$ws.Range("K2").Formula = '=COUNTIF(F2:F200,"TP")'
$ws.Range("L2").Formula = '=COUNTIF(G2:G200,"FP")'
$ws.Range("M2").Formula = '=COUNTIF(H2:H200,"TN")'
$ws.Range("N2").Formula = '=COUNTIF(I2:I200,"FN")'

# --- Match the authored selection state (I2:I68, active cell I2). ---
$ws.Range("I2:I68").Select()

$wb.Save()
